# CPI corrected, new Plots
# Update CSCC values (column C, rows 2-6) with corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 491.9042826293455
$ws.Range("C3").Value = 138.3149247278208
$ws.Range("C4").Value = 25.80198226248687
$ws.Range("C5").Value = 56.19539220537131
$ws.Range("C6").Value = -32.76315078674113
